# إضافة حدث جديد في Card2 by admin at 2025-12-11 08:56:55
# Fill in row 18's previously-blank B:K columns with "nan" placeholder text,
# then append a brand-new event row (row 19) to the Card2 worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card2")

# --- Normalize row 18 (B18:K18 were blank placeholders -> "nan") ---
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")
foreach ($col in $cols) {
    $ws.Range($col + "18").Value = "nan"
}

# --- Append new event row 19 ---
# Column A holds the card number, stored as text (matches existing rows).
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "2"
$ws.Range("A19").ClearFormats()

# Columns B:K are left blank for the new row, same as row 18 originally was.

# Column L holds the date, stored as text (must not auto-convert to a date serial).
$ws.Range("L19").NumberFormat = "@"
$ws.Range("L19").Value = "11/12/2025"
$ws.Range("L19").ClearFormats()

# Column M: Event description.
$ws.Range("M19").Value = "قطع سير 700"

# Column N: Correction description.
$ws.Range("N19").Value = "تم تغير سير مشرشر  700مجموعه الكلندرات"

# Column O: Serviced by.
$ws.Range("O19").Value = "م/محمد عبدالله/ايهاب"

"Row 19 added to Card2"
